$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings; some target values look like plain numbers
# (e.g. "1.01"). Excel would silently convert those to real numbers on
# assignment, so format those specific cells as Text first to keep them as
# strings, matching the source data (cells whose value is not ambiguous are
# left with default/general formatting, unchanged from the original).
$ws.Range("D5:D6").NumberFormat = "@"
$ws.Range("D8:D13").NumberFormat = "@"
$ws.Range("D15:D16").NumberFormat = "@"
$ws.Range("D19:D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '60.751.48'
$ws.Cells.Item(2, 5).Value = '  +6.48%  '
$ws.Cells.Item(3, 4).Value = '3.320.90'
$ws.Cells.Item(3, 5).Value = '  +2.62%  '
$ws.Cells.Item(4, 5).Value = '  +0.96%  '
$ws.Cells.Item(5, 4).Value = '407.37'
$ws.Cells.Item(5, 5).Value = '  +3.05%  '
$ws.Cells.Item(6, 4).Value = '110.28'
$ws.Cells.Item(6, 5).Value = '  +3.14%  '
$ws.Cells.Item(7, 2).Value = 'LidoStakedEther'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Cells.Item(7, 4).Value = '3.289.53'
$ws.Cells.Item(7, 5).Value = '  +1.65%  '
$ws.Cells.Item(8, 2).Value = 'XRP'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Cells.Item(8, 4).Value = '0.557'
$ws.Cells.Item(8, 5).Value = '  -2.20%  '
$ws.Cells.Item(9, 2).Value = 'USDC'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Cells.Item(9, 4).Value = '1.00'
$ws.Cells.Item(9, 5).Value = '  -0.02%  '
$ws.Cells.Item(10, 2).Value = 'Cardano'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Cells.Item(10, 4).Value = '0.620'
$ws.Cells.Item(10, 5).Value = '  +0.96%  '
$ws.Cells.Item(11, 2).Value = 'Dogecoin'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(11, 4).Value = '0.108'
$ws.Cells.Item(11, 5).Value = '  +10.82%  '
$ws.Cells.Item(12, 2).Value = 'Avalanche'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(12, 4).Value = '38.45'
$ws.Cells.Item(12, 5).Value = '  -1.11%  '
$ws.Cells.Item(13, 2).Value = 'TRON'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(13, 4).Value = '0.142'
$ws.Cells.Item(13, 5).Value = '  +0.20%  '
$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).Value = '4.015.51'
$ws.Cells.Item(14, 5).Value = '  +6.98%  '
$ws.Cells.Item(15, 4).Value = '18.79'
$ws.Cells.Item(15, 5).Value = '  -0.59%  '
$ws.Cells.Item(16, 2).Value = 'Polkadot'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(16, 4).Value = '7.97'
$ws.Cells.Item(16, 5).Value = '  -1.43%  '
$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).Value = '3.505.88'
$ws.Cells.Item(17, 5).Value = '  +8.03%  '
$ws.Cells.Item(18, 2).Value = 'WrappedBTC'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(18, 4).Value = '61.458.32'
$ws.Cells.Item(18, 5).Value = '  +7.99%  '
$ws.Cells.Item(19, 2).Value = 'Polygon'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(19, 4).Value = '0.994'
$ws.Cells.Item(19, 5).Value = '  -3.52%  '
$ws.Cells.Item(20, 2).Value = 'Uniswap'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(20, 4).Value = '10.30'
$ws.Cells.Item(20, 5).Value = '  -6.17%  '
$ws.Cells.Item(21, 2).Value = 'ShibaInu'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(21, 4).Value = '0.0000108'
$ws.Cells.Item(21, 5).Value = '  +2.35%  '
$ws.Cells.Item(22, 2).Value = 'ImmutableX'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(22, 4).Value = '3.22'
$ws.Cells.Item(22, 5).Value = '  -3.25%  '
$ws.Cells.Item(23, 4).Value = '289.14'
$ws.Cells.Item(23, 5).Value = '  -2.35%  '
$ws.Cells.Item(24, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(24, 4).Value = '11.82'
$ws.Cells.Item(24, 5).Value = '  -9.02%  '
$ws.Cells.Item(25, 2).Value = 'Litecoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(25, 4).Value = '73.66'
$ws.Cells.Item(25, 5).Value = '  -0.25%  '
$ws.Cells.Item(26, 2).Value = 'PancakeSwap'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(26, 4).Value = '3.19'
$ws.Cells.Item(26, 5).Value = '  +1.57%  '
$ws.Cells.Item(27, 2).Value = 'EthereumClassic'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(27, 4).Value = '28.24'
$ws.Cells.Item(27, 5).Value = '  +1.55%  '
$ws.Cells.Item(28, 2).Value = 'LEO'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(28, 4).Value = '4.51'
$ws.Cells.Item(28, 5).Value = '  +2.89%  '
$ws.Cells.Item(29, 4).Value = '0.158'
$ws.Cells.Item(29, 5).Value = '  -6.46%  '
$ws.Cells.Item(30, 4).Value = '7.12'
$ws.Cells.Item(30, 5).Value = '  -1.38%  '
$ws.Cells.Item(31, 2).Value = 'Filecoin'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(31, 4).Value = '7.29'
$ws.Cells.Item(31, 5).Value = '  -4.93%  '
$ws.Cells.Item(32, 2).Value = 'Hedera'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(32, 4).Value = '0.108'
$ws.Cells.Item(32, 5).Value = '  +0.08%  '
$ws.Cells.Item(33, 4).Value = '0.998'
$ws.Cells.Item(33, 5).Value = '  -0.19%  '
$ws.Cells.Item(34, 2).Value = 'Cosmos'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(34, 4).Value = '11.04'
$ws.Cells.Item(34, 5).Value = '  +0.06%  '
$ws.Cells.Item(35, 4).Value = '37.07'
$ws.Cells.Item(35, 5).Value = '  -0.60%  '
$ws.Cells.Item(36, 2).Value = 'Toncoin'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(36, 4).Value = '2.36'
$ws.Cells.Item(36, 5).Value = '  +11.58%  '
$ws.Cells.Item(37, 2).Value = 'VeChain'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(37, 4).Value = '0.0470'
$ws.Cells.Item(37, 5).Value = '  -2.92%  '
$ws.Cells.Item(38, 2).Value = 'OKB'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(38, 4).Value = '52.12'
$ws.Cells.Item(38, 5).Value = '  +0.89%  '
$ws.Cells.Item(39, 4).Value = '1.01'
$ws.Cells.Item(39, 5).Value = '  +0.59%  '
$ws.Cells.Item(40, 2).Value = 'Stacks'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(40, 4).Value = '2.98'
$ws.Cells.Item(40, 5).Value = '  -1.72%  '
$ws.Cells.Item(41, 2).Value = 'EnergySwap'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(41, 4).Value = '27.36'
$ws.Cells.Item(41, 5).Value = '  +25.19%  '
$ws.Cells.Item(42, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(42, 4).Value = '3.24'
$ws.Cells.Item(42, 5).Value = '  -8.49%  '
$ws.Cells.Item(43, 2).Value = 'Monero'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(43, 4).Value = '137.38'
$ws.Cells.Item(43, 5).Value = '  +2.24%  '
$ws.Cells.Item(44, 2).Value = 'Stellar'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(44, 4).Value = '0.118'
$ws.Cells.Item(44, 5).Value = '  -1.39%  '
$ws.Cells.Item(45, 2).Value = 'ARBITRUM'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(45, 4).Value = '1.85'
$ws.Cells.Item(45, 5).Value = '  -2.04%  '
$ws.Cells.Item(46, 2).Value = 'TheGraph'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Cells.Item(46, 4).Value = '0.275'
$ws.Cells.Item(46, 5).Value = '  -1.95%  '
$ws.Cells.Item(47, 2).Value = 'NEARProtocol'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(47, 4).Value = '3.72'
$ws.Cells.Item(47, 5).Value = '  -5.89%  '
$ws.Cells.Item(48, 2).Value = 'Celestia'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Cells.Item(48, 4).Value = '15.89'
$ws.Cells.Item(48, 5).Value = '  -5.74%  '
$ws.Cells.Item(49, 2).Value = 'WEMIXToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(49, 4).Value = '2.25'
$ws.Cells.Item(49, 5).Value = '  +7.69%  '
$ws.Cells.Item(50, 2).Value = 'Maker'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(50, 4).Value = '2.122.75'
$ws.Cells.Item(50, 5).Value = '  -1.16%  '
$ws.Cells.Item(51, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(51, 4).Value = '2.37'
$ws.Cells.Item(51, 5).Value = '  +1.24%  '
